# Applies the diff: for rows 2-20 on the "Specification" sheet, insert the
# value "Application" into column C, shifting the existing C/D/E values
# rightward into D/E/F respectively (column F in this range was always
# empty, so its previous value is simply discarded). Columns A, B, G, H, I
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 20; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 6).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
    $ws.Cells.Item($r, 3).Value = "Application"
}
